$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Valor") updates.
# Values ending in "%" would otherwise be auto-parsed by Excel into a
# percentage number, losing the literal "x.xx %" text - a leading
# apostrophe forces literal text entry, matching how these were typed
# originally (shared-string cells, not numbers).
$ws.Range("C4").Value = "457055.01 N"

$ws.Range("C10").Value = "'6.89 %"
$ws.Range("C11").Value = "'4.27 %"
$ws.Range("C12").Value = "'3.70 %"
$ws.Range("C13").Value = "'9.75 %"
$ws.Range("C14").Value = "'19.92 %"
$ws.Range("C15").Value = "'6.78 %"

$ws.Range("C18").Value = "9106.99 N (Ref T0 req)"
$ws.Range("C19").Value = "'8.90 %"
$ws.Range("C20").Value = "'4.54 %"

$ws.Range("C30").Value = "'9.21 %"
$ws.Range("C31").Value = "'7.09 %"
$ws.Range("C32").Value = "15.76 º"
$ws.Range("C33").Value = "19.33 º"
$ws.Range("C34").Value = "58.93 º"
$ws.Range("C35").Value = "'88.97 %"
$ws.Range("C36").Value = "30.10 m"

# Column D ("Status") updates
$ws.Range("D20").Value = "NOK"
$ws.Range("D31").Value = "OK"
$ws.Range("D32").Value = "OK"
$ws.Range("D33").Value = "OK"
$ws.Range("D35").Value = "OK"
